$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 701.8
$ws.Range("I41").Value = 992.7778
$ws.Range("J41").Value = 265.33334
$ws.Range("K41").Value = 992.7778
$ws.Range("L41").Value = 265.33334
$ws.Range("M41").Value = -552.7778
$ws.Range("N41").Value = -1145.33334

$ws.Range("H80").Value = 461.45456
$ws.Range("I80").Value = 547.1667
$ws.Range("J80").Value = 358.6
$ws.Range("K80").Value = 1641.5001
$ws.Range("L80").Value = 1075.8
$ws.Range("M80").Value = -643.5001
$ws.Range("N80").Value = -3071.8

$ws.Range("H83").Value = 461.45456
$ws.Range("I83").Value = 547.1667
$ws.Range("J83").Value = 358.6
$ws.Range("K83").Value = 4924.5003
$ws.Range("L83").Value = 3227.4
$ws.Range("M83").Value = 67.4997000000003
$ws.Range("N83").Value = -13211.4

$ws.Range("H132").Value = 1450.4125
$ws.Range("I132").Value = 1586.4714
$ws.Range("J132").Value = 498
$ws.Range("K132").Value = 4759.414199999999
$ws.Range("L132").Value = 1494
$ws.Range("M132").Value = -2229.414199999999
$ws.Range("N132").Value = -6554

$ws.Range("H137").Value = 97242.484
$ws.Range("I137").Value = 1520.0892
$ws.Range("J137").Value = 480132.06
$ws.Range("K137").Value = 4560.267599999999
$ws.Range("L137").Value = 1440396.18
$ws.Range("M137").Value = -2010.267599999999
$ws.Range("N137").Value = -1445496.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5729.778
$ws.Range("I32").Value = 3381.803
$ws.Range("J32").Value = 31557.5
$ws.Range("K32").Value = 3381.803
$ws.Range("L32").Value = 31557.5
$ws.Range("M32").Value = -3094.803
$ws.Range("N32").Value = -32131.5

$ws.Range("H74").Value = 138387.73
$ws.Range("I74").Value = 197214.48
$ws.Range("J74").Value = 2016.6364
$ws.Range("K74").Value = 197214.48
$ws.Range("L74").Value = 2016.6364
$ws.Range("M74").Value = -196340.48
$ws.Range("N74").Value = -3764.6364

$ws.Range("H77").Value = 138387.73
$ws.Range("I77").Value = 197214.48
$ws.Range("J77").Value = 2016.6364
$ws.Range("K77").Value = 986072.4
$ws.Range("L77").Value = 10083.182
$ws.Range("M77").Value = -981704.4
$ws.Range("N77").Value = -18819.182

$ws.Range("H132").Value = 8060.51
$ws.Range("I132").Value = 13035.5
$ws.Range("J132").Value = 2886.52
$ws.Range("K132").Value = 39106.5
$ws.Range("L132").Value = 8659.559999999999
$ws.Range("M132").Value = -36576.5
$ws.Range("N132").Value = -13719.56

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1707.3846
$ws.Range("I99").Value = 1889.7
$ws.Range("J99").Value = 1099.6666
$ws.Range("K99").Value = 1889.7
$ws.Range("L99").Value = 1099.6666
$ws.Range("M99").Value = -391.7
$ws.Range("N99").Value = -4095.6666

$ws.Range("H102").Value = 9988.375
$ws.Range("I102").Value = 4969.2
$ws.Range("J102").Value = 18353.666
$ws.Range("K102").Value = 4969.2
$ws.Range("L102").Value = 18353.666
$ws.Range("M102").Value = -1724.2
$ws.Range("N102").Value = -24843.666

$ws.Range("H138").Value = 90000
$ws.Range("I138").Value = 90000
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 90000
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -84860
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 8411.866
$ws.Range("I25").Value = 2066.3333
$ws.Range("J25").Value = 9998.25
$ws.Range("K25").Value = 2066.3333
$ws.Range("L25").Value = 9998.25
$ws.Range("M25").Value = -1892.3333
$ws.Range("N25").Value = -10346.25

$ws.Range("H31").Value = 2373.3735
$ws.Range("I31").Value = 1978.6578
$ws.Range("J31").Value = 2706.689
$ws.Range("K31").Value = 1978.6578
$ws.Range("L31").Value = 2706.689
$ws.Range("M31").Value = -1683.6578
$ws.Range("N31").Value = -3296.689

$ws.Range("H34").Value = 2373.3735
$ws.Range("I34").Value = 1978.6578
$ws.Range("J34").Value = 2706.689
$ws.Range("K34").Value = 1978.6578
$ws.Range("L34").Value = 2706.689
$ws.Range("M34").Value = -1776.6578
$ws.Range("N34").Value = -3110.689

$ws.Range("H47").Value = 19032
$ws.Range("I47").Value = 18064
$ws.Range("J47").Value = 20000
$ws.Range("K47").Value = 18064
$ws.Range("L47").Value = 20000
$ws.Range("M47").Value = -17498
$ws.Range("N47").Value = -21132

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()

$ws.Range("H55").Value = 3358333.2
$ws.Range("I55").Value = 5020000
$ws.Range("J55").Value = 35000
$ws.Range("K55").Value = 5020000
$ws.Range("L55").Value = 35000
$ws.Range("M55").Value = -5019685
$ws.Range("N55").Value = -35630

$ws.Range("H58").Value = 1791.6478
$ws.Range("I58").Value = 1563.8276
$ws.Range("J58").Value = 2808.077
$ws.Range("K58").Value = 1563.8276
$ws.Range("L58").Value = 2808.077
$ws.Range("M58").Value = -1360.8276
$ws.Range("N58").Value = -3214.077

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H99").Value = 3475.2
$ws.Range("I99").Value = 3632.9375
$ws.Range("J99").Value = 2844.25
$ws.Range("K99").Value = 3632.9375
$ws.Range("L99").Value = 2844.25
$ws.Range("M99").Value = -2134.9375
$ws.Range("N99").Value = -5840.25

$ws.Range("H126").Value = 3475.2
$ws.Range("I126").Value = 3632.9375
$ws.Range("J126").Value = 2844.25
$ws.Range("K126").Value = 10898.8125
$ws.Range("L126").Value = 8532.75
$ws.Range("M126").Value = -8428.8125
$ws.Range("N126").Value = -13472.75

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("N135").Value = 0
$ws.Range("L135").ClearContents()

$ws.Range("H136").Value = 1791.6478
$ws.Range("I136").Value = 1563.8276
$ws.Range("J136").Value = 2808.077
$ws.Range("K136").Value = 4691.4828
$ws.Range("L136").Value = 8424.231
$ws.Range("M136").Value = -2141.4828
$ws.Range("N136").Value = -13524.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1388
$ws.Range("I7").Value = 82
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 246
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = -134
$ws.Range("N7").Value = -12224

$ws.Range("H107").Value = 1291.8695
$ws.Range("I107").Value = 306
$ws.Range("J107").Value = 1925.6428
$ws.Range("K107").Value = 918
$ws.Range("L107").Value = 5776.928400000001
$ws.Range("M107").Value = 1002
$ws.Range("N107").Value = -9616.928400000001

$ws.Range("H122").Value = 1946.1428
$ws.Range("I122").Value = 925.5
$ws.Range("J122").Value = 2354.4
$ws.Range("K122").Value = 8329.5
$ws.Range("L122").Value = 21189.6
$ws.Range("M122").Value = -5879.5
$ws.Range("N122").Value = -26089.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 62500
$ws.Range("I26").Value = 65000
$ws.Range("J26").Value = 60000
$ws.Range("K26").Value = 65000
$ws.Range("L26").Value = 60000
$ws.Range("M26").Value = -64720
$ws.Range("N26").Value = -60560

$ws.Range("H50").Value = 62500
$ws.Range("I50").Value = 65000
$ws.Range("J50").Value = 60000
$ws.Range("K50").Value = 65000
$ws.Range("L50").Value = 60000
$ws.Range("M50").Value = -64502
$ws.Range("N50").Value = -60996

$ws.Range("H113").Value = 83335336
$ws.Range("I113").Value = 90911064
$ws.Range("J113").Value = 2313
$ws.Range("K113").Value = 90911064
$ws.Range("L113").Value = 2313
$ws.Range("M113").Value = -90908894
$ws.Range("N113").Value = -6653

$ws.Range("H132").Value = 5062.273
$ws.Range("I132").Value = 5308.6665
$ws.Range("J132").Value = 3953.5
$ws.Range("K132").Value = 15925.9995
$ws.Range("L132").Value = 11860.5
$ws.Range("M132").Value = -13395.9995
$ws.Range("N132").Value = -16920.5

$ws.Range("H135").Value = 99535.27
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 99535.27
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 99535.27
$ws.Range("N135").Value = -109675.27
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 6455.222
$ws.Range("I17").Value = 5555
$ws.Range("J17").Value = 6712.4287
$ws.Range("K17").Value = 5555
$ws.Range("L17").Value = 6712.4287
$ws.Range("M17").Value = -5385
$ws.Range("N17").Value = -7052.4287

$ws.Range("H30").Value = 18755.076
$ws.Range("I30").Value = 20222.666
$ws.Range("J30").Value = 17497.143
$ws.Range("K30").Value = 20222.666
$ws.Range("L30").Value = 17497.143
$ws.Range("M30").Value = -20114.666
$ws.Range("N30").Value = -17713.143

$ws.Range("H31").Value = 14334.591
$ws.Range("I31").Value = 11101.615
$ws.Range("J31").Value = 19004.445
$ws.Range("K31").Value = 11101.615
$ws.Range("L31").Value = 19004.445
$ws.Range("M31").Value = -10853.615
$ws.Range("N31").Value = -19500.445

$ws.Range("H46").Value = 2754.724
$ws.Range("I46").Value = 1128.4
$ws.Range("J46").Value = 4497.2144
$ws.Range("K46").Value = 1128.4
$ws.Range("L46").Value = 4497.2144
$ws.Range("M46").Value = -940.4000000000001
$ws.Range("N46").Value = -4873.2144

$ws.Range("H61").Value = 34484412
$ws.Range("I61").Value = 38462810
$ws.Range("J61").Value = 4965.6665
$ws.Range("K61").Value = 38462810
$ws.Range("L61").Value = 4965.6665
$ws.Range("M61").Value = -38462608
$ws.Range("N61").Value = -5369.6665

$ws.Range("H93").Value = 940.9231
$ws.Range("I93").Value = 940.9231
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 940.9231
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 307.0769

$ws.Range("H113").Value = 34484412
$ws.Range("I113").Value = 38462810
$ws.Range("J113").Value = 4965.6665
$ws.Range("K113").Value = 38462810
$ws.Range("L113").Value = 4965.6665
$ws.Range("M113").Value = -38460640
$ws.Range("N113").Value = -9305.666499999999

$ws.Range("H132").Value = 2848.54
$ws.Range("I132").Value = 2474.081
$ws.Range("J132").Value = 3914.3076
$ws.Range("K132").Value = 7422.243
$ws.Range("L132").Value = 11742.9228
$ws.Range("M132").Value = -4892.243
$ws.Range("N132").Value = -16802.9228

$ws.Range("H136").Value = 19934.736
$ws.Range("I136").Value = 2629.2307
$ws.Range("J136").Value = 68142.92999999999
$ws.Range("K136").Value = 7887.6921
$ws.Range("L136").Value = 204428.79
$ws.Range("M136").Value = -5337.6921
$ws.Range("N136").Value = -209528.79

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3788.303
$ws.Range("I122").Value = 3509.7407
$ws.Range("J122").Value = 5041.8335
$ws.Range("K122").Value = 10529.2221
$ws.Range("L122").Value = 15125.5005
$ws.Range("M122").Value = -8079.222099999999
$ws.Range("N122").Value = -20025.5005

$ws.Range("H126").Value = 2452.8235
$ws.Range("I126").Value = 2707.3076
$ws.Range("J126").Value = 1625.75
$ws.Range("K126").Value = 8121.9228
$ws.Range("L126").Value = 4877.25
$ws.Range("M126").Value = -5651.9228
$ws.Range("N126").Value = -9817.25

$ws.Range("H132").Value = 2323.2273
$ws.Range("I132").Value = 2395.2812
$ws.Range("J132").Value = 2131.0833
$ws.Range("K132").Value = 7185.8436
$ws.Range("L132").Value = 6393.249899999999
$ws.Range("M132").Value = -4655.8436
$ws.Range("N132").Value = -11453.2499

$ws.Range("H136").Value = 4241.4683
$ws.Range("I136").Value = 4536.25
$ws.Range("J136").Value = 3276.7273
$ws.Range("K136").Value = 13608.75
$ws.Range("L136").Value = 9830.1819
$ws.Range("M136").Value = -10830.1819
$ws.Range("N136").Value = -14930.1819
